$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.266.44"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "3.962.04"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'611.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").Value = "'170.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("D7").Value = "3.960.32"
$ws.Range("E7").Value = "  +3.44%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "'0.172"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.13%  "
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("E13").Value = "  +5.73%  "
$ws.Range("D14").Value = "'38.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").Value = "4.629.56"
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("D16").Value = "3.937.84"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "70.217.34"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "'7.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'17.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "'502.68"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.35%  "
$ws.Range("D23").Value = "'0.745"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "'0.0000170"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.64%  "
$ws.Range("D25").Value = "'85.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").Value = "'10.29"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'3.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "4.122.00"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("D33").Value = "'7.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'32.48"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "3.936.68"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").Value = "'3.28"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.94%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("E43").Value = "  +5.20%  "
$ws.Range("D44").Value = "'440.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "'0.000281"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +25.24%  "
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").Value = "'40.73"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.89%  "
$ws.Range("D51").Value = "'143.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.33%  "
